$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Migrate Interfaces")

# The original data contained two erroneous duplicate rows caused by a
# parse error when the port-channel interfaces were imported:
#   - Row 10 duplicated the "Ethernet1/7" / port-channel 7 entry with
#     incorrect VPC id / VLAN / description data (belongs further down).
#   - Row 13 duplicated the "Ethernet1/21" entry with incorrect VPC id /
#     switchport mode / description data.
# Deleting these two rows lets the remaining rows (and their correct
# per-row banding style) shift up into place, matching the corrected
# export.

$ws.Rows.Item(13).Delete()
$ws.Rows.Item(10).Delete()
